$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = [double]"7.327471962526033e-15"
$ws.Range("C2").Value = [double]"2.052465086777033e-09"
$ws.Range("D2").Value = [double]"3.223369029078222"
$ws.Range("E2").Value = [double]"13.86384647080068"
$ws.Range("G2").Value = [double]"17.08721550193138"

# Row 3
$ws.Range("B3").Value = [double]"0.1169995834814548"
$ws.Range("C3").Value = [double]"0.3048912486333797"
$ws.Range("D3").Value = [double]"3.223369029078222"
$ws.Range("E3").Value = [double]"13.86384647080068"
$ws.Range("G3").Value = [double]"17.50910633199374"

# Row 4
$ws.Range("B4").Value = [double]"0.6545652718822623"
$ws.Range("C4").Value = [double]"1.626987699542094"
$ws.Range("D4").Value = [double]"0.7210945179870265"
$ws.Range("E4").Value = [double]"0.5333859586016987"
$ws.Range("G4").Value = [double]"3.536033448013082"
